$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E8 from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Select E8 as the active cell
$ws.Range("E8").Select()
